# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Row -> new F-column value, for each of the two sheets sharing this data.
$updates = @{
    2  = 20
    4  = 271
    5  = 46
    6  = 550
    7  = 54
    8  = 2006
    10 = 97
    11 = 4313
    13 = 281
    15 = 3
    16 = 108
    17 = 23
    18 = 15
    19 = 64
    20 = 3079
    21 = 64
    22 = 450
    25 = 72
    26 = 76
    29 = 51
    30 = 197
    31 = 9
    32 = 479
    33 = 1701
    34 = 253
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
